# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    4  = 68
    5  = 531
    6  = 7312
    7  = 209
    8  = 177
    9  = 1062
    10 = 468
    11 = 12
    12 = 155
    14 = 661
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
